$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.782.13"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.414.13"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.89"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.80"
$ws.Range("E6").Value = "  +4.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.427.95"
$ws.Range("E9").Value = "  +1.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  +2.79%  "

$ws.Range("E11").Value = "  +0.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  +3.97%  "

$ws.Range("E13").Value = "  +3.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.59"
$ws.Range("E14").Value = "  +3.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  +4.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.834.02"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.369.90"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.421.18"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").Value = "  +1.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.84"
$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("E22").Value = "  +1.53%  "

$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  +12.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.14"
$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "629.45"
$ws.Range("E27").Value = "  +11.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.43"
$ws.Range("E28").Value = "  +2.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0956"
$ws.Range("E29").Value = "  +4.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.528.12"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.09"
$ws.Range("E31").Value = "  +1.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  +6.32%  "

$ws.Range("E33").Value = "  +0.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.134"
$ws.Range("E34").Value = "  +1.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("E35").Value = "  +4.98%  "

$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  +2.75%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.372"
$ws.Range("E38").Value = "  +1.56%  "

$ws.Range("E39").Value = "  -0.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.39"
$ws.Range("E40").Value = "  +5.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.44"
$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").Value = "  +9.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.72"
$ws.Range("E43").Value = "  +4.22%  "

$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.08"
$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("E46").Value = "  -2.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.91"
$ws.Range("E47").Value = "  +0.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.56"
$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.98"
$ws.Range("E49").Value = "  +4.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.597"
$ws.Range("E50").Value = "  +2.05%  "

$ws.Range("E51").Value = "  +2.53%  "
